# Update "想去人数" (column F) values on the "展览" and "全部类型" worksheets
# to reflect newly generated data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F, shared by both sheets.
$updates = @{
    3  = 2909
    5  = 151
    7  = 1573
    11 = 1300
    13 = 419
    18 = 86
    20 = 2951
    23 = 75
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
